$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1808510638297872
$ws.Range("C2").Value = 0.5638297872340425
$ws.Range("J2").Value = 0.02482269503546099
$ws.Range("P2").Value = 0.1347517730496454
$ws.Range("S2").Value = 0.09574468085106383
$ws.Range("B3").Value = 0.02762430939226519
$ws.Range("C3").Value = 0.04419889502762431
$ws.Range("J3").Value = 0.02762430939226519
$ws.Range("P3").Value = 0.7292817679558011
$ws.Range("S3").Value = 0.1712707182320442
$ws.Range("J4").Value = 0.03571428571428571
$ws.Range("P4").Value = 0.8214285714285714
$ws.Range("P5").Value = 0.8
$ws.Range("S5").Value = 0.2
$ws.Range("B6").Value = 0.03125
$ws.Range("F6").Value = 0.046875
$ws.Range("J6").Value = 0.2135416666666667
$ws.Range("O6").Value = 0.04166666666666666
$ws.Range("Q6").Value = 0.1614583333333333
$ws.Range("R6").Value = 0.078125
$ws.Range("S6").Value = 0.4270833333333333
$ws.Range("B7").Value = 0.1235955056179775
$ws.Range("F7").Value = 0.02247191011235955
$ws.Range("J7").Value = 0.1235955056179775
$ws.Range("O7").Value = 0.02247191011235955
$ws.Range("Q7").Value = 0.1629213483146068
$ws.Range("R7").Value = 0.07865168539325842
$ws.Range("S7").Value = 0.4662921348314606
$ws.Range("B8").Value = 0.09641873278236915
$ws.Range("D8").Value = 0.01377410468319559
$ws.Range("F8").Value = 0.04958677685950413
$ws.Range("J8").Value = 0.1267217630853995
$ws.Range("O8").Value = 0.005509641873278237
$ws.Range("Q8").Value = 0.1129476584022039
$ws.Range("R8").Value = 0.1101928374655647
$ws.Range("S8").Value = 0.4848484848484849
$ws.Range("B9").Value = 0.1103896103896104
$ws.Range("D9").Value = 0.01948051948051948
$ws.Range("E9").Value = 0.006493506493506494
$ws.Range("F9").Value = 0.07792207792207792
$ws.Range("J9").Value = 0.1103896103896104
$ws.Range("O9").Value = 0.006493506493506494
$ws.Range("Q9").Value = 0.1363636363636364
$ws.Range("R9").Value = 0.06493506493506493
$ws.Range("S9").Value = 0.4675324675324675
$ws.Range("B10").Value = 0.1229946524064171
$ws.Range("D10").Value = 0.01693404634581105
$ws.Range("E10").Value = 0.0035650623885918
$ws.Range("F10").Value = 0.08645276292335116
$ws.Range("J10").Value = 0.1140819964349376
$ws.Range("O10").Value = 0.008021390374331552
$ws.Range("Q10").Value = 0.1827094474153298
$ws.Range("R10").Value = 0.07397504456327986
$ws.Range("S10").Value = 0.3912655971479501
$ws.Range("G11").Value = 0.1584158415841584
$ws.Range("J11").Value = 0.0891089108910891
$ws.Range("K11").Value = 0.2079207920792079
$ws.Range("L11").Value = 0.5148514851485149
$ws.Range("S11").Value = 0.0297029702970297
$ws.Range("G12").Value = 0.7345679012345679
$ws.Range("J12").Value = 0.191358024691358
$ws.Range("L12").Value = 0.02469135802469136
$ws.Range("S12").Value = 0.04938271604938271
$ws.Range("G13").Value = 0.5833333333333334
$ws.Range("J13").Value = 0.3333333333333333
$ws.Range("S13").Value = 0.08333333333333333
$ws.Range("F15").Value = 0.04046242774566474
$ws.Range("H15").Value = 0.1213872832369942
$ws.Range("I15").Value = 0.05780346820809248
$ws.Range("J15").Value = 0.3526011560693642
$ws.Range("K15").Value = 0.1329479768786127
$ws.Range("M15").Value = 0.0115606936416185
$ws.Range("O15").Value = 0.05202312138728324
$ws.Range("S15").Value = 0.2312138728323699
$ws.Range("F16").Value = 0.005405405405405406
$ws.Range("H16").Value = 0.1567567567567568
$ws.Range("I16").Value = 0.04864864864864865
$ws.Range("J16").Value = 0.3945945945945946
$ws.Range("K16").Value = 0.1351351351351351
$ws.Range("M16").Value = 0.03243243243243243
$ws.Range("N16").Value = 0.005405405405405406
$ws.Range("O16").Value = 0.05945945945945946
$ws.Range("S16").Value = 0.1621621621621622
$ws.Range("F17").Value = 0.01829268292682927
$ws.Range("H17").Value = 0.1341463414634146
$ws.Range("I17").Value = 0.09451219512195122
$ws.Range("J17").Value = 0.4817073170731707
$ws.Range("K17").Value = 0.08231707317073171
$ws.Range("M17").Value = 0.01219512195121951
$ws.Range("O17").Value = 0.06097560975609756
$ws.Range("S17").Value = 0.1158536585365854
$ws.Range("H18").Value = 0.1715976331360947
$ws.Range("I18").Value = 0.05917159763313609
$ws.Range("J18").Value = 0.4497041420118343
$ws.Range("K18").Value = 0.0650887573964497
$ws.Range("M18").Value = 0.02958579881656805
$ws.Range("O18").Value = 0.03550295857988166
$ws.Range("S18").Value = 0.1893491124260355
$ws.Range("F19").Value = 0.01331114808652246
$ws.Range("H19").Value = 0.2038269550748752
$ws.Range("I19").Value = 0.07986688851913477
$ws.Range("J19").Value = 0.3627287853577371
$ws.Range("K19").Value = 0.1331114808652246
$ws.Range("M19").Value = 0.01580698835274542
$ws.Range("N19").Value = 0.001663893510815308
$ws.Range("O19").Value = 0.06655574043261231
$ws.Range("S19").Value = 0.1231281198003328
